$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: update S18, clear V18 and W18
$ws.Range("S18").Value = 58430
$ws.Range("V18").ClearContents()
$ws.Range("W18").ClearContents()

# Row 19: add S19
$ws.Range("S19").Value = 56357

# Row 20: add S20
$ws.Range("S20").Value = 55191

# Row 21: add S21
$ws.Range("S21").Value = 52357

# Row 22: add S22
$ws.Range("S22").Value = 50614

# Row 23: add S23, V23, W23
$ws.Range("S23").Value = 47951
$ws.Range("V23").Value = 3675.599486111112
$ws.Range("W23").Value = 6.987983166606228
